$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Detail rows: stock-rule recalculation ---

# Row 3
$ws.Range("L3").Value = 0
$ws.Range("P3").Value = 29
$ws.Range("Q3").Value = 0
$ws.Range("U3").Value = 0
$ws.Rows.Item(3).Hidden = $true

# Row 6
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("P6").Value = 16
$ws.Range("Q6").Value = 0
$ws.Range("U6").Value = 0
$ws.Rows.Item(6).Hidden = $true

# Row 8
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("P8").Value = 35
$ws.Range("Q8").Value = 0
$ws.Range("U8").Value = 0
$ws.Rows.Item(8).Hidden = $true

# Row 11
$ws.Range("L11").Value = 0

# Row 14
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("P14").Value = 11
$ws.Range("Q14").Value = 0
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 2
$ws.Range("U14").Value = 0
$ws.Rows.Item(14).Hidden = $true

# Row 17
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("P17").Value = 6
$ws.Range("Q17").Value = 0
$ws.Range("U17").Value = 0
$ws.Rows.Item(17).Hidden = $true

# Row 19
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("P19").Value = 23
$ws.Range("Q19").Value = 0
$ws.Range("U19").Value = 0
$ws.Rows.Item(19).Hidden = $true

# Row 20
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("P20").Value = 6
$ws.Range("Q20").Value = 0
$ws.Range("U20").Value = 0
$ws.Rows.Item(20).Hidden = $true

# Row 21
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("P21").Value = 12
$ws.Range("Q21").Value = 0
$ws.Range("U21").Value = 0
$ws.Rows.Item(21).Hidden = $true

# Row 23
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("P23").Value = 16
$ws.Range("Q23").Value = 0
$ws.Range("U23").Value = 0
$ws.Rows.Item(23).Hidden = $true

# Row 24
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("P24").Value = 17
$ws.Range("Q24").Value = 0
$ws.Range("S24").Value = 1
$ws.Range("T24").Value = 1
$ws.Range("U24").Value = 0
$ws.Rows.Item(24).Hidden = $true

# Row 25
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("P25").Value = 30
$ws.Range("Q25").Value = 0
$ws.Range("S25").Value = 3
$ws.Range("T25").Value = 3
$ws.Range("U25").Value = 0
$ws.Rows.Item(25).Hidden = $true

# Row 26
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("P26").Value = 90
$ws.Range("Q26").Value = 0
$ws.Range("S26").Value = 3
$ws.Range("T26").Value = 3
$ws.Range("U26").Value = 0
$ws.Rows.Item(26).Hidden = $true

# Row 27
$ws.Range("L27").Value = 0

# Row 28
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("P28").Value = 10
$ws.Range("Q28").Value = 0
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = 1
$ws.Range("U28").Value = 0
$ws.Rows.Item(28).Hidden = $true

# Row 29
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = 0
$ws.Range("P29").Value = 24
$ws.Range("Q29").Value = 0
$ws.Range("U29").Value = 0
$ws.Rows.Item(29).Hidden = $true

# Row 31
$ws.Range("M31").Value = 16.65
$ws.Range("N31").Value = 9.99
$ws.Range("P31").Value = 3
$ws.Range("Q31").Value = 0
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("U31").Value = 1

# Row 32
$ws.Range("M32").Value = 22.38
$ws.Range("N32").Value = 13.43
$ws.Range("P32").Value = 1
$ws.Range("Q32").Value = 1
$ws.Range("U32").Value = 1

# Row 33
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("P33").Value = 19
$ws.Range("Q33").Value = 0
$ws.Range("U33").Value = 0
$ws.Rows.Item(33).Hidden = $true

# Row 34
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("P34").Value = 11
$ws.Range("Q34").Value = 0
$ws.Range("U34").Value = 0
$ws.Rows.Item(34).Hidden = $true

# --- Summary metrics (rows 38-49) ---

$ws.Range("C38").Value = 42

# C40 contains a literal text value with a euro sign; avoid the automatic
# locale-aware number/currency parsing that happens on direct .Value string
# assignment by writing it as a text formula and then collapsing the
# formula down to its static result via copy / paste-values.
$c40 = $ws.Range("C40")
$c40.Formula = '=T("506.59€")'
$c40.Copy()
$c40.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("C49").Value = 0
